# Fix traj read real robot data: update columns B-F for rows 1-6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.1619081784284657
$ws.Range("C1").Value = 0.1945208102290912
$ws.Range("D1").Value = -0.672860002190442
$ws.Range("E1").Value = 0.7034158221868738
$ws.Range("F1").Value = 1.570796280873402

$ws.Range("B2").Value = 0.2074942785077096
$ws.Range("C2").Value = 0.1934383557349265
$ws.Range("D2").Value = -0.6758762528771421
$ws.Range("E2").Value = 0.7014819866207826
$ws.Range("F2").Value = 1.570796286790826

$ws.Range("B3").Value = 0.4117502296915038
$ws.Range("C3").Value = 0.1885882419516801
$ws.Range("D3").Value = -0.6893910556777709
$ws.Range("E3").Value = 0.6928171211837666
$ws.Range("F3").Value = 1.570796313304808

$ws.Range("B4").Value = 0.6991189600253006
$ws.Range("C4").Value = 0.1817645923503434
$ws.Range("D4").Value = -0.7084051000618862
$ws.Range("E4").Value = 0.6806264781953137
$ws.Range("F4").Value = 1.570796350607462

$ws.Range("B5").Value = 0.9033749112090953
$ws.Range("C5").Value = 0.176914478567097
$ws.Range("D5").Value = -0.7219199028625151
$ws.Range("E5").Value = 0.6719616127582977
$ws.Range("F5").Value = 1.570796377121443

$ws.Range("B6").Value = 0.9489610112883384
$ws.Range("C6").Value = 0.1758320240729323
$ws.Range("D6").Value = -0.7249361535492151
$ws.Range("E6").Value = 0.6700277771922065
$ws.Range("F6").Value = 1.570796383038867
